# "Generate Report for Archive"
# - Localization status moves from "Ready for handoff" to "In Translation"
#   for every file/language row (Overview!E2:F4 and the per-language
#   Status column, zh-cn!C2:C4 / de-de!C2:C4).
# - The Status columns are narrower now that the new text is shorter than
#   the old one, so re-size them to match the refreshed content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns -------
$overview = $wb.Worksheets.Item("Overview")
$overviewStatusCols = @("E", "F")
foreach ($col in $overviewStatusCols) {
    for ($row = 2; $row -le 4; $row++) {
        $cell = $overview.Range($col + $row)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Per-language sheets: Status column (col C) --------------------------
$languageSheets = @("zh-cn", "de-de")
foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 4; $row++) {
        $cell = $ws.Range("C" + $row)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Re-fit the Status columns now that the text is shorter --------------
# ColumnWidth is quantized to whole pixels by the host, so 12.5 is the
# nearest reachable value to the narrower width the shorter text needs.
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # Overview!E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # Overview!F (de-de)

foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(3).ColumnWidth = $newColumnWidth     # Status column
}
